$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Original row 17 (GABALDON 2021 - Camanjac ES, Dumaguete City) is removed entirely;
# every subsequent row shifts up by one.
$ws.Rows.Item(17).Delete()

# Original row 29 (GABALDON 2024 - South City Central School, Toledo City) -
# now sitting at row 28 after the first deletion - is also removed, so the
# data block ends at row 27.
$ws.Rows.Item(28).Delete()

# Column R (18th column) narrows from 43 to 27 characters.
$ws.Columns.Item(18).ColumnWidth = 157/6
